# Update the LR-pair data rows with recomputed statistics (per Dr Hou advice).
# Each of the 4 sending clusters now has values for all 4 target clusters (16 rows total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgals1"
$ws.Range("C2").Value = "Ptprc"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 25.02674933333333
$ws.Range("H2").Value = 75.080248
$ws.Range("I2").Value = 0.05787790829091637
$ws.Range("J2").Value = 0.05787790829091637
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 205.313027
$ws.Range("N2").Value = 615.9390810000001
$ws.Range("O2").Value = 0.435242422384838
$ws.Range("P2").Value = 0.435242422384838
$ws.Range("Q2").Value = 5138.317661596899
$ws.Range("R2").Value = 46244.8589543721
$ws.Range("S2").Value = 0.02519092100710594
$ws.Range("T2").Value = 0.02519092100710594

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgals1"
$ws.Range("C3").Value = "Ptprc"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 25.02674933333333
$ws.Range("H3").Value = 75.080248
$ws.Range("I3").Value = 0.05787790829091637
$ws.Range("J3").Value = 0.05787790829091637
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.186821
$ws.Range("N3").Value = 0.5604629999999999
$ws.Range("O3").Value = 0.0003960412341120362
$ws.Range("P3").Value = 0.0003960412341120362
$ws.Range("Q3").Value = 4.675522337202666
$ws.Range("R3").Value = 42.07970103482399
$ws.Range("S3").Value = 0.00002292203822735777
$ws.Range("T3").Value = 0.00002292203822735777

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lgals1"
$ws.Range("C4").Value = "Ptprc"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 25.02674933333333
$ws.Range("H4").Value = 75.080248
$ws.Range("I4").Value = 0.05787790829091637
$ws.Range("J4").Value = 0.05787790829091637
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 266.1765593333333
$ws.Range("N4").Value = 798.529678
$ws.Range("O4").Value = 0.5642668278730386
$ws.Range("P4").Value = 0.5642668278730386
$ws.Range("Q4").Value = 6661.534028844459
$ws.Range("R4").Value = 59953.80625960014
$ws.Range("S4").Value = 0.03265858371524202
$ws.Range("T4").Value = 0.03265858371524202

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lgals1"
$ws.Range("C5").Value = "Ptprc"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 25.02674933333333
$ws.Range("H5").Value = 75.080248
$ws.Range("I5").Value = 0.05787790829091637
$ws.Range("J5").Value = 0.05787790829091637
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.044676
$ws.Range("N5").Value = 0.134028
$ws.Range("O5").Value = 0.00009470850801135487
$ws.Range("P5").Value = 0.00009470850801135488
$ws.Range("Q5").Value = 1.118095053216
$ws.Range("R5").Value = 10.062855478944
$ws.Range("S5").Value = 0.000005481530341050716
$ws.Range("T5").Value = 0.000005481530341050716

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lgals1"
$ws.Range("C6").Value = "Ptprc"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 138.9376323333333
$ws.Range("H6").Value = 416.812897
$ws.Range("I6").Value = 0.3213129853678317
$ws.Range("J6").Value = 0.3213129853678316
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 205.313027
$ws.Range("N6").Value = 615.9390810000001
$ws.Range("O6").Value = 0.435242422384838
$ws.Range("P6").Value = 0.435242422384838
$ws.Range("Q6").Value = 28525.70585856975
$ws.Range("R6").Value = 256731.3527271277
$ws.Range("S6").Value = 0.1398490420951991
$ws.Range("T6").Value = 0.139849042095199

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lgals1"
$ws.Range("C7").Value = "Ptprc"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 138.9376323333333
$ws.Range("H7").Value = 416.812897
$ws.Range("I7").Value = 0.3213129853678317
$ws.Range("J7").Value = 0.3213129853678316
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.186821
$ws.Range("N7").Value = 0.5604629999999999
$ws.Range("O7").Value = 0.0003960412341120362
$ws.Range("P7").Value = 0.0003960412341120362
$ws.Range("Q7").Value = 25.95646741014567
$ws.Range("R7").Value = 233.608206691311
$ws.Range("S7").Value = 0.0001272531912612987
$ws.Range("T7").Value = 0.0001272531912612987

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Lgals1"
$ws.Range("C8").Value = "Ptprc"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 138.9376323333333
$ws.Range("H8").Value = 416.812897
$ws.Range("I8").Value = 0.3213129853678317
$ws.Range("J8").Value = 0.3213129853678316
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 266.1765593333333
$ws.Range("N8").Value = 798.529678
$ws.Range("O8").Value = 0.5642668278730386
$ws.Range("P8").Value = 0.5642668278730386
$ws.Range("Q8").Value = 36981.94093640635
$ws.Range("R8").Value = 332837.4684276571
$ws.Range("S8").Value = 0.1813062590079225
$ws.Range("T8").Value = 0.1813062590079224

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Lgals1"
$ws.Range("C9").Value = "Ptprc"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 138.9376323333333
$ws.Range("H9").Value = 416.812897
$ws.Range("I9").Value = 0.3213129853678317
$ws.Range("J9").Value = 0.3213129853678316
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.044676
$ws.Range("N9").Value = 0.134028
$ws.Range("O9").Value = 0.00009470850801135487
$ws.Range("P9").Value = 0.00009470850801135488
$ws.Range("Q9").Value = 6.207177662124001
$ws.Range("R9").Value = 55.86459895911601
$ws.Range("S9").Value = 0.00003043107344886164
$ws.Range("T9").Value = 0.00003043107344886164

# Row 10: M2 -> ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Lgals1"
$ws.Range("C10").Value = "Ptprc"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 98.29468800000001
$ws.Range("H10").Value = 294.884064
$ws.Range("I10").Value = 0.2273204107243322
$ws.Range("J10").Value = 0.2273204107243321
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 205.313027
$ws.Range("N10").Value = 615.9390810000001
$ws.Range("O10").Value = 0.435242422384838
$ws.Range("P10").Value = 0.435242422384838
$ws.Range("Q10").Value = 20181.17993130058
$ws.Range("R10").Value = 181630.6193817052
$ws.Range("S10").Value = 0.09893948622117464
$ws.Range("T10").Value = 0.09893948622117461

# Row 11: M2 -> FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Lgals1"
$ws.Range("C11").Value = "Ptprc"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 98.29468800000001
$ws.Range("H11").Value = 294.884064
$ws.Range("I11").Value = 0.2273204107243322
$ws.Range("J11").Value = 0.2273204107243321
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.186821
$ws.Range("N11").Value = 0.5604629999999999
$ws.Range("O11").Value = 0.0003960412341120362
$ws.Range("P11").Value = 0.0003960412341120362
$ws.Range("Q11").Value = 18.363511906848
$ws.Range("R11").Value = 165.271607161632
$ws.Range("S11").Value = 0.00009002825600211945
$ws.Range("T11").Value = 0.00009002825600211942

# Row 12: M2 -> M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Lgals1"
$ws.Range("C12").Value = "Ptprc"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 98.29468800000001
$ws.Range("H12").Value = 294.884064
$ws.Range("I12").Value = 0.2273204107243322
$ws.Range("J12").Value = 0.2273204107243321
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 266.1765593333333
$ws.Range("N12").Value = 798.529678
$ws.Range("O12").Value = 0.5642668278730386
$ws.Range("P12").Value = 0.5642668278730386
$ws.Range("Q12").Value = 26163.74185258349
$ws.Range("R12").Value = 235473.6766732514
$ws.Range("S12").Value = 0.1282693670702152
$ws.Range("T12").Value = 0.1282693670702151

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Lgals1"
$ws.Range("C13").Value = "Ptprc"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 98.29468800000001
$ws.Range("H13").Value = 294.884064
$ws.Range("I13").Value = 0.2273204107243322
$ws.Range("J13").Value = 0.2273204107243321
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.044676
$ws.Range("N13").Value = 0.134028
$ws.Range("O13").Value = 0.00009470850801135487
$ws.Range("P13").Value = 0.00009470850801135488
$ws.Range("Q13").Value = 4.391413481088001
$ws.Range("R13").Value = 39.52272132979201
$ws.Range("S13").Value = 0.00002152917694022989
$ws.Range("T13").Value = 0.00002152917694022989

# Row 14: sCs -> ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Lgals1"
$ws.Range("C14").Value = "Ptprc"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 170.1468356666666
$ws.Range("H14").Value = 510.440507
$ws.Range("I14").Value = 0.3934886956169198
$ws.Range("J14").Value = 0.3934886956169198
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 205.313027
$ws.Range("N14").Value = 615.9390810000001
$ws.Range("O14").Value = 0.435242422384838
$ws.Range("P14").Value = 0.435242422384838
$ws.Range("Q14").Value = 34933.3618651949
$ws.Range("R14").Value = 314400.2567867541
$ws.Range("S14").Value = 0.1712629730613584
$ws.Range("T14").Value = 0.1712629730613584

# Row 15: sCs -> FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Lgals1"
$ws.Range("C15").Value = "Ptprc"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 170.1468356666666
$ws.Range("H15").Value = 510.440507
$ws.Range("I15").Value = 0.3934886956169198
$ws.Range("J15").Value = 0.3934886956169198
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.186821
$ws.Range("N15").Value = 0.5604629999999999
$ws.Range("O15").Value = 0.0003960412341120362
$ws.Range("P15").Value = 0.0003960412341120362
$ws.Range("Q15").Value = 31.78700198608233
$ws.Range("R15").Value = 286.083017874741
$ws.Range("S15").Value = 0.0001558377486212603
$ws.Range("T15").Value = 0.0001558377486212603

# Row 16: sCs -> M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Lgals1"
$ws.Range("C16").Value = "Ptprc"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 170.1468356666666
$ws.Range("H16").Value = 510.440507
$ws.Range("I16").Value = 0.3934886956169198
$ws.Range("J16").Value = 0.3934886956169198
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 266.1765593333333
$ws.Range("N16").Value = 798.529678
$ws.Range("O16").Value = 0.5642668278730386
$ws.Range("P16").Value = 0.5642668278730386
$ws.Range("Q16").Value = 45289.09929920741
$ws.Range("R16").Value = 407601.8936928667
$ws.Range("S16").Value = 0.222032618079659
$ws.Range("T16").Value = 0.222032618079659

# Row 17: sCs -> sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Lgals1"
$ws.Range("C17").Value = "Ptprc"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 170.1468356666666
$ws.Range("H17").Value = 510.440507
$ws.Range("I17").Value = 0.3934886956169198
$ws.Range("J17").Value = 0.3934886956169198
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.044676
$ws.Range("N17").Value = 0.134028
$ws.Range("O17").Value = 0.00009470850801135487
$ws.Range("P17").Value = 0.00009470850801135488
$ws.Range("Q17").Value = 7.601480030243999
$ws.Range("R17").Value = 68.413320272196
$ws.Range("S17").Value = 0.00003726672728121263
$ws.Range("T17").Value = 0.00003726672728121263
